$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells for columns I and J
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Match the header formatting used by the rest of row 1 (e.g. H1: bold, bordered, centered)
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats

# Data values for column I (rows 2-11)
$colI = @(9, 8, 9, 9, 9, 9, 9, 9, 9, 9)
# Data values for column J (rows 2-11)
$colJ = @(9, 9, 9, 9, 9, 9, 9, 9, 9, 9)

for ($i = 0; $i -lt 10; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $colI[$i]
    $ws.Cells.Item($row, 10).Value = $colJ[$i]
}
